$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 was mapped to field_wbddh_terms_of_use; now it is mapped to the
# License field instead, with a header value in column A.
$ws.Range("B9").Value = "field_license_wbddh"
$ws.Range("A9").Value = "License"

# Update the active selection to A10 (as captured in the sheet view).
$ws.Range("A10").Select()
